$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1597.1052  # H19: was 1542.2
$ws.Cells.Item(19, 10).Value = 526.8570999999999  # J19: was 523.375
$ws.Cells.Item(19, 12).Value = 526.8570999999999  # L19: was 523.375
$ws.Cells.Item(19, 14).Value = -876.8570999999999  # N19: was -873.375

$ws.Cells.Item(64, 8).Value = 5840  # H64: was 4449.8
$ws.Cells.Item(64, 9).Value = 3600  # I64: was 3375
$ws.Cells.Item(64, 10).Value = 7333.3335  # J64: was 5166.3335
$ws.Cells.Item(64, 11).Value = 3600  # K64: was 3375
$ws.Cells.Item(64, 12).Value = 7333.3335  # L64: was 5166.3335
$ws.Cells.Item(64, 13).Value = -3352  # M64: was -3127
$ws.Cells.Item(64, 14).Value = -7829.3335  # N64: was -5662.3335

$ws.Cells.Item(67, 8).Value = 5840  # H67: was 4449.8
$ws.Cells.Item(67, 9).Value = 3600  # I67: was 3375
$ws.Cells.Item(67, 10).Value = 7333.3335  # J67: was 5166.3335
$ws.Cells.Item(67, 11).Value = 3600  # K67: was 3375
$ws.Cells.Item(67, 12).Value = 7333.3335  # L67: was 5166.3335
$ws.Cells.Item(67, 13).Value = -2742  # M67: was -2517
$ws.Cells.Item(67, 14).Value = -9049.333500000001  # N67: was -6882.3335

$ws.Cells.Item(94, 8).Value = 743.73334  # H94: was 768
$ws.Cells.Item(94, 9).Value = 755.7857  # I94: was 768
$ws.Cells.Item(94, 10).Value = 575  # J94: was 0
$ws.Cells.Item(94, 11).Value = 755.7857  # K94: was 768
$ws.Cells.Item(94, 12).Value = 575  # L94: was 0
$ws.Cells.Item(94, 13).Value = -304.7857  # M94: was -317
$ws.Cells.Item(94, 14).Value = -1477  # N94: was None

$ws.Cells.Item(116, 8).Value = 5758.1177  # H116: was 6706.5
$ws.Cells.Item(116, 9).Value = 3845.2856  # I116: was 4752.5
$ws.Cells.Item(116, 10).Value = 7097.1  # J116: was 7097.3
$ws.Cells.Item(116, 11).Value = 3845.2856  # K116: was 4752.5
$ws.Cells.Item(116, 12).Value = 7097.1  # L116: was 7097.3
$ws.Cells.Item(116, 13).Value = -403.2856000000002  # M116: was -1310.5
$ws.Cells.Item(116, 14).Value = -13981.1  # N116: was -13981.3

$ws.Cells.Item(125, 8).Value = 5036  # H125: was 2500
$ws.Cells.Item(125, 9).Value = 0  # I125: was 2500
$ws.Cells.Item(125, 10).Value = 5036  # J125: was 0
$ws.Cells.Item(125, 11).Value = 0  # K125: was 22500
$ws.Cells.Item(125, 12).ClearContents()  # L125: was 0
$ws.Cells.Item(125, 13).Value = 45324  # M125: was -20040
$ws.Cells.Item(125, 14).Value = -50244  # N125: was None

$ws.Cells.Item(132, 8).Value = 1843.8334  # H132: was 1850.973
$ws.Cells.Item(132, 9).Value = 1899.4828  # I132: was 1954.9656
$ws.Cells.Item(132, 10).Value = 1613.2858  # J132: was 1474
$ws.Cells.Item(132, 11).Value = 5698.4484  # K132: was 5864.8968
$ws.Cells.Item(132, 12).Value = 4839.857400000001  # L132: was 4422
$ws.Cells.Item(132, 13).Value = -3168.4484  # M132: was -3334.8968
$ws.Cells.Item(132, 14).Value = -9899.857400000001  # N132: was -9482

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 9481.143  # H45: was 5899.5864
$ws.Cells.Item(45, 9).Value = 11434.454  # I45: was 6735.9165
$ws.Cells.Item(45, 10).Value = 2319  # J45: was 1885.2
$ws.Cells.Item(45, 11).Value = 11434.454  # K45: was 6735.9165
$ws.Cells.Item(45, 12).Value = 2319  # L45: was 1885.2
$ws.Cells.Item(45, 13).Value = -11057.454  # M45: was -6358.9165
$ws.Cells.Item(45, 14).Value = -3073  # N45: was -2639.2

$ws.Cells.Item(97, 8).Value = 1056.0278  # H97: was 1453
$ws.Cells.Item(97, 9).Value = 970.24243  # I97: was 1577.7222
$ws.Cells.Item(97, 10).Value = 1999.6666  # J97: was 1004
$ws.Cells.Item(97, 11).Value = 970.24243  # K97: was 1577.7222
$ws.Cells.Item(97, 12).Value = 1999.6666  # L97: was 1004
$ws.Cells.Item(97, 13).Value = -474.24243  # M97: was -1081.7222
$ws.Cells.Item(97, 14).Value = -2991.6666  # N97: was -1996

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 873.5217  # H94: was 1777.4445
$ws.Cells.Item(94, 9).Value = 898.3143  # I94: was 1086.4615
$ws.Cells.Item(94, 10).Value = 794.63635  # J94: was 3574
$ws.Cells.Item(94, 11).Value = 898.3143  # K94: was 1086.4615
$ws.Cells.Item(94, 12).Value = 794.63635  # L94: was 3574
$ws.Cells.Item(94, 13).Value = -447.3143  # M94: was -635.4614999999999
$ws.Cells.Item(94, 14).Value = -1696.63635  # N94: was -4476

$ws.Cells.Item(134, 8).Value = 6905.569  # H134: was 6791.981
$ws.Cells.Item(134, 9).Value = 3185.244  # I134: was 3133.1904
$ws.Cells.Item(134, 11).Value = 9555.732  # K134: was 9399.5712
$ws.Cells.Item(134, 13).Value = -7020.732  # M134: was -6864.5712

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4749.9546  # H31: was 4650.6665
$ws.Cells.Item(31, 9).Value = 4149.4  # I31: was 3710.1667
$ws.Cells.Item(31, 10).Value = 5250.4165  # J31: was 5591.1665
$ws.Cells.Item(31, 11).Value = 4149.4  # K31: was 3710.1667
$ws.Cells.Item(31, 12).Value = 5250.4165  # L31: was 5591.1665
$ws.Cells.Item(31, 13).Value = -3854.4  # M31: was -3415.1667
$ws.Cells.Item(31, 14).Value = -5840.4165  # N31: was -6181.1665

$ws.Cells.Item(34, 8).Value = 4749.9546  # H34: was 4650.6665
$ws.Cells.Item(34, 9).Value = 4149.4  # I34: was 3710.1667
$ws.Cells.Item(34, 10).Value = 5250.4165  # J34: was 5591.1665
$ws.Cells.Item(34, 11).Value = 4149.4  # K34: was 3710.1667
$ws.Cells.Item(34, 12).Value = 5250.4165  # L34: was 5591.1665
$ws.Cells.Item(34, 13).Value = -3947.4  # M34: was -3508.1667
$ws.Cells.Item(34, 14).Value = -5654.4165  # N34: was -5995.1665

$ws.Cells.Item(58, 8).Value = 3052.3076  # H58: was 4789.643
$ws.Cells.Item(58, 9).Value = 1502.0435  # I58: was 2081.2307
$ws.Cells.Item(58, 10).Value = 14937.667  # J58: was 39999
$ws.Cells.Item(58, 11).Value = 1502.0435  # K58: was 2081.2307
$ws.Cells.Item(58, 12).Value = 14937.667  # L58: was 39999
$ws.Cells.Item(58, 13).Value = -1299.0435  # M58: was -1878.2307
$ws.Cells.Item(58, 14).Value = -15343.667  # N58: was -40405

$ws.Cells.Item(136, 8).Value = 3052.3076  # H136: was 4789.643
$ws.Cells.Item(136, 9).Value = 1502.0435  # I136: was 2081.2307
$ws.Cells.Item(136, 10).Value = 14937.667  # J136: was 39999
$ws.Cells.Item(136, 11).Value = 4506.1305  # K136: was 6243.6921
$ws.Cells.Item(136, 12).Value = 44813.001  # L136: was 119997
$ws.Cells.Item(136, 13).Value = -1956.1305  # M136: was -3693.6921
$ws.Cells.Item(136, 14).Value = -49913.001  # N136: was -125097

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 2187.625  # H80: was 2583.8333
$ws.Cells.Item(80, 9).Value = 1874.5  # I80: was 2750
$ws.Cells.Item(80, 11).Value = 5623.5  # K80: was 8250
$ws.Cells.Item(80, 13).Value = -4687.5  # M80: was -7314

$ws.Cells.Item(83, 8).Value = 2187.625  # H83: was 2583.8333
$ws.Cells.Item(83, 9).Value = 1874.5  # I83: was 2750
$ws.Cells.Item(83, 11).Value = 16870.5  # K83: was 24750
$ws.Cells.Item(83, 13).Value = -12190.5  # M83: was -20070

$ws.Cells.Item(132, 8).Value = 3032542  # H132: was 100000000
$ws.Cells.Item(132, 9).Value = 2074.35  # I132: was 0
$ws.Cells.Item(132, 10).Value = 7694799.5  # J132: was 100000000
$ws.Cells.Item(132, 11).Value = 18669.15  # K132: was 0
$ws.Cells.Item(132, 12).Value = 69253195.5  # L132: was 900000000
$ws.Cells.Item(132, 13).Value = -16139.15  # M132: was None
$ws.Cells.Item(132, 14).Value = -69258255.5  # N132: was -900005060

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 22146  # H26: was 22574.572
$ws.Cells.Item(26, 9).Value = 12000  # I26: was 0
$ws.Cells.Item(26, 10).Value = 23837  # J26: was 22574.572
$ws.Cells.Item(26, 11).Value = 12000  # K26: was 0
$ws.Cells.Item(26, 12).Value = 23837  # L26: was 22574.572
$ws.Cells.Item(26, 13).Value = -11720  # M26: was None
$ws.Cells.Item(26, 14).Value = -24397  # N26: was -23134.572

$ws.Cells.Item(50, 8).Value = 22146  # H50: was 22574.572
$ws.Cells.Item(50, 9).Value = 12000  # I50: was 0
$ws.Cells.Item(50, 10).Value = 23837  # J50: was 22574.572
$ws.Cells.Item(50, 11).Value = 12000  # K50: was 0
$ws.Cells.Item(50, 12).Value = 23837  # L50: was 22574.572
$ws.Cells.Item(50, 13).Value = -11502  # M50: was None
$ws.Cells.Item(50, 14).Value = -24833  # N50: was -23570.572

$ws.Cells.Item(122, 8).Value = 5484.543  # H122: was 5682.951
$ws.Cells.Item(122, 9).Value = 4152.0186  # I122: was 4546.5093
$ws.Cells.Item(122, 10).Value = 8149.593  # J122: was 7759.8965
$ws.Cells.Item(122, 11).Value = 12456.0558  # K122: was 13639.5279
$ws.Cells.Item(122, 12).Value = 24448.779  # L122: was 23279.6895
$ws.Cells.Item(122, 13).Value = -10006.0558  # M122: was -11189.5279
$ws.Cells.Item(122, 14).Value = -29348.779  # N122: was -28179.6895

$ws.Cells.Item(132, 8).Value = 15386.259  # H132: was 14860.75
$ws.Cells.Item(132, 9).Value = 11465.091  # I132: was 10587.667
$ws.Cells.Item(132, 10).Value = 32639.4  # J132: was 40499.25
$ws.Cells.Item(132, 11).Value = 34395.273  # K132: was 31763.001
$ws.Cells.Item(132, 12).Value = 97918.20000000001  # L132: was 121497.75
$ws.Cells.Item(132, 13).Value = -31865.273  # M132: was -29233.001
$ws.Cells.Item(132, 14).Value = -102978.2  # N132: was -126557.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3694.6338  # H132: was 3794.3623
$ws.Cells.Item(132, 9).Value = 2870.8215  # I132: was 2986.6
$ws.Cells.Item(132, 10).Value = 6770.2  # J132: was 6967.7144
$ws.Cells.Item(132, 11).Value = 8612.4645  # K132: was 8959.799999999999
$ws.Cells.Item(132, 12).Value = 20310.6  # L132: was 20903.1432
$ws.Cells.Item(132, 13).Value = -6082.4645  # M132: was -6429.799999999999
$ws.Cells.Item(132, 14).Value = -25370.6  # N132: was -25963.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 15833.333  # H51: was 20000
$ws.Cells.Item(51, 9).Value = 15833.333  # I51: was 20000
$ws.Cells.Item(51, 11).Value = 15833.333  # K51: was 20000
$ws.Cells.Item(51, 13).Value = -15323.333  # M51: was -19490

$ws.Cells.Item(52, 8).Value = 18680.334  # H52: was 16020.5
$ws.Cells.Item(52, 9).Value = 18680.334  # I52: was 16020.5
$ws.Cells.Item(52, 11).Value = 18680.334  # K52: was 16020.5
$ws.Cells.Item(52, 13).Value = -18454.334  # M52: was -15794.5

$ws.Cells.Item(132, 8).Value = 5726.69  # H132: was 12595.258
$ws.Cells.Item(132, 9).Value = 4707.5435  # I132: was 7565.761
$ws.Cells.Item(132, 10).Value = 17446.875  # J132: was 25555.885
$ws.Cells.Item(132, 11).Value = 14122.6305  # K132: was 22697.283
$ws.Cells.Item(132, 12).Value = 52340.625  # L132: was 76667.655
$ws.Cells.Item(132, 13).Value = -11592.6305  # M132: was -20167.283
$ws.Cells.Item(132, 14).Value = -57400.625  # N132: was -81727.655
